$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Select()

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Validate_validPinCodeAtRegistrationPage"
$ws.Cells.Item(14, 3).Value = "Functional"
$ws.Cells.Item(14, 4).Value = "verify & validate that application is checking valid and invalid pin code format in registration page."
$ws.Cells.Item(14, 4).Style = $ws.Cells.Item(13,4).Style

$ws.Range("C22").Select()
